# Update "想去人数" (number interested) values on the "展览" and "全部类型" sheets
# to match the freshly generated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new numeric value
$updates = @{
    "F2" = 2343
    "F3" = 1855
    "F6" = 1105
    "F7" = 48
    "F8" = 5950
    "F9" = 98
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
